$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '62.086.43'
$ws.Cells.Item(2, 5).Value = '  +1.93%  '

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.434.59'
$ws.Cells.Item(3, 5).Value = '  +2.25%  '

$ws.Cells.Item(4, 5).Value = '  -0.09%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '408.64'
$ws.Cells.Item(5, 5).Value = '  +1.18%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '129.44'
$ws.Cells.Item(6, 5).Value = '  -2.76%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.627'
$ws.Cells.Item(7, 5).Value = '  +6.58%  '

$ws.Cells.Item(8, 5).Value = '  -0.08%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.748'
$ws.Cells.Item(9, 5).Value = '  +12.16%  '

$ws.Cells.Item(10, 5).Value = '  +18.77%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '43.46'
$ws.Cells.Item(11, 5).Value = '  +3.26%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.141'
$ws.Cells.Item(12, 5).Value = '  -0.25%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '8.72'
$ws.Cells.Item(13, 5).Value = '  +5.23%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '20.24'
$ws.Cells.Item(14, 5).Value = '  +3.92%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.0000197'
$ws.Cells.Item(15, 5).Value = '  +56.62%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '3.471.88'
$ws.Cells.Item(16, 5).Value = '  +2.85%  '

$ws.Cells.Item(17, 2).Value = 'WrappedBTC'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '62.086.74'
$ws.Cells.Item(17, 5).Value = '  +1.81%  '

$ws.Cells.Item(18, 2).Value = 'Polygon'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '1.04'
$ws.Cells.Item(18, 5).Value = '  +2.96%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '11.43'
$ws.Cells.Item(19, 5).Value = '  +3.68%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '371.48'
$ws.Cells.Item(20, 5).Value = '  +22.54%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '86.60'
$ws.Cells.Item(21, 5).Value = '  +3.81%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '3.18'
$ws.Cells.Item(22, 5).Value = '  -0.84%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '13.25'
$ws.Cells.Item(23, 5).Value = '  +4.50%  '

$ws.Cells.Item(24, 5).Value = '  +2.60%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '31.65'
$ws.Cells.Item(25, 5).Value = '  +8.21%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '4.81'
$ws.Cells.Item(26, 5).Value = '  +0.56%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '8.38'
$ws.Cells.Item(27, 5).Value = '  +1.15%  '

$ws.Cells.Item(28, 5).Value = '  +2.59%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.72'
$ws.Cells.Item(29, 5).Value = '  +9.96%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '44.15'
$ws.Cells.Item(30, 5).Value = '  +7.23%  '

$ws.Cells.Item(31, 5).Value = '  -0.61%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.116'
$ws.Cells.Item(32, 5).Value = '  +0.26%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '11.78'
$ws.Cells.Item(33, 5).Value = '  +4.54%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.999'
$ws.Cells.Item(34, 5).Value = '  +0.01%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.0491'
$ws.Cells.Item(35, 5).Value = '  +3.00%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '52.01'
$ws.Cells.Item(36, 5).Value = '  +0.21%  '

$ws.Cells.Item(37, 5).Value = '  -0.06%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '3.36'
$ws.Cells.Item(38, 5).Value = '  -1.35%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.91'
$ws.Cells.Item(39, 5).Value = '  +0.87%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.314'
$ws.Cells.Item(40, 5).Value = '  +9.54%  '

$ws.Cells.Item(41, 2).Value = 'Stellar'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.131'
$ws.Cells.Item(41, 5).Value = '  +6.39%  '

$ws.Cells.Item(42, 2).Value = 'Monero'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '143.02'
$ws.Cells.Item(42, 5).Value = '  +4.31%  '

$ws.Cells.Item(43, 5).Value = '  +0.34%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '4.00'
$ws.Cells.Item(44, 5).Value = '  +0.96%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '16.65'
$ws.Cells.Item(45, 5).Value = '  +0.30%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.33'
$ws.Cells.Item(46, 5).Value = '  +4.39%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '21.73'
$ws.Cells.Item(47, 5).Value = '  +1.79%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.116.65'
$ws.Cells.Item(48, 5).Value = '  +0.21%  '

$ws.Cells.Item(49, 5).Value = '  -0.29%  '

$ws.Cells.Item(50, 5).Value = '  +2.91%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.0363'
$ws.Cells.Item(51, 5).Value = '  +7.69%  '
